$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 331.9091  # H2
$ws.Cells.Item(2, 9).Value = 315  # I2
$ws.Cells.Item(2, 10).Value = 341.57144  # J2
$ws.Cells.Item(2, 11).Value = 315  # K2
$ws.Cells.Item(2, 12).Value = 341.57144  # L2
$ws.Cells.Item(2, 13).Value = -202  # M2
$ws.Cells.Item(2, 14).Value = -567.5714399999999  # N2
$ws.Cells.Item(42, 8).Value = 2349.8  # H42
$ws.Cells.Item(42, 9).Value = 5150  # I42
$ws.Cells.Item(42, 11).Value = 15450  # K42
$ws.Cells.Item(42, 13).Value = -15220  # M42
$ws.Cells.Item(70, 8).Value = 3268.7727  # H70
$ws.Cells.Item(70, 9).Value = 2463.6365  # I70
$ws.Cells.Item(70, 10).Value = 4073.9092  # J70
$ws.Cells.Item(70, 11).Value = 7390.9095  # K70
$ws.Cells.Item(70, 12).Value = 12221.7276  # L70
$ws.Cells.Item(70, 13).Value = -7120.9095  # M70
$ws.Cells.Item(70, 14).Value = -12761.7276  # N70
$ws.Cells.Item(73, 8).Value = 3268.7727  # H73
$ws.Cells.Item(73, 9).Value = 2463.6365  # I73
$ws.Cells.Item(73, 10).Value = 4073.9092  # J73
$ws.Cells.Item(73, 11).Value = 7390.9095  # K73
$ws.Cells.Item(73, 12).Value = 12221.7276  # L73
$ws.Cells.Item(73, 13).Value = -6454.9095  # M73
$ws.Cells.Item(73, 14).Value = -14093.7276  # N73
$ws.Cells.Item(112, 8).Value = 0  # H112
$ws.Cells.Item(112, 9).Value = 0  # I112
$ws.Cells.Item(112, 10).Value = 0  # J112
$ws.Cells.Item(112, 11).Value = 0  # K112
$ws.Cells.Item(112, 12).ClearContents()  # L112
$ws.Cells.Item(112, 13).ClearContents()  # M112
$ws.Cells.Item(112, 14).Value = 0  # N112
$ws.Cells.Item(133, 8).Value = 65000  # H133
$ws.Cells.Item(133, 10).Value = 65000  # J133
$ws.Cells.Item(133, 12).Value = 65000  # L133
$ws.Cells.Item(133, 14).Value = -75120  # N133
$ws.Cells.Item(135, 8).Value = 192.5  # H135
$ws.Cells.Item(135, 9).Value = 175  # I135
$ws.Cells.Item(135, 10).Value = 245  # J135
$ws.Cells.Item(135, 11).Value = 1575  # K135
$ws.Cells.Item(135, 12).Value = 2205  # L135
$ws.Cells.Item(135, 13).Value = 960  # M135
$ws.Cells.Item(135, 14).Value = -7275  # N135
$ws.Cells.Item(138, 8).Value = 2457.875  # H138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 1851.25  # H4
$ws.Cells.Item(4, 9).Value = 579.6  # I4
$ws.Cells.Item(4, 10).Value = 3970.6667  # J4
$ws.Cells.Item(4, 11).Value = 579.6  # K4
$ws.Cells.Item(4, 12).Value = 3970.6667  # L4
$ws.Cells.Item(4, 13).Value = -463.6  # M4
$ws.Cells.Item(4, 14).Value = -4202.6667  # N4
$ws.Cells.Item(29, 8).Value = 5639.75  # H29
$ws.Cells.Item(29, 10).Value = 7386.3335  # J29
$ws.Cells.Item(29, 12).Value = 7386.3335  # L29
$ws.Cells.Item(29, 14).Value = -8002.3335  # N29
$ws.Cells.Item(30, 8).Value = 962.25  # H30
$ws.Cells.Item(30, 9).Value = 1146.3334  # I30
$ws.Cells.Item(30, 11).Value = 1146.3334  # K30
$ws.Cells.Item(30, 13).Value = -996.3334  # M30
$ws.Cells.Item(38, 8).Value = 3949999.8  # H38
$ws.Cells.Item(38, 10).Value = 4933333  # J38
$ws.Cells.Item(38, 12).Value = 4933333  # L38
$ws.Cells.Item(38, 14).Value = -4934267  # N38
$ws.Cells.Item(122, 8).Value = 8978.6  # H122
$ws.Cells.Item(122, 9).Value = 8723.5  # I122
$ws.Cells.Item(122, 10).Value = 9999  # J122
$ws.Cells.Item(122, 11).Value = 26170.5  # K122
$ws.Cells.Item(122, 12).Value = 29997  # L122
$ws.Cells.Item(122, 13).Value = -23720.5  # M122
$ws.Cells.Item(122, 14).Value = -34897  # N122
$ws.Cells.Item(132, 8).Value = 1710.7  # H132
$ws.Cells.Item(132, 9).Value = 1763.375  # I132
$ws.Cells.Item(132, 11).Value = 5290.125  # K132
$ws.Cells.Item(132, 13).Value = -2760.125  # M132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 105.8  # H11
$ws.Cells.Item(11, 9).Value = 50.5  # I11
$ws.Cells.Item(11, 10).Value = 142.66667  # J11
$ws.Cells.Item(11, 11).Value = 50.5  # K11
$ws.Cells.Item(11, 12).Value = 142.66667  # L11
$ws.Cells.Item(11, 13).Value = 89.5  # M11
$ws.Cells.Item(11, 14).Value = -422.66667  # N11
$ws.Cells.Item(33, 8).Value = 13673.667  # H33
$ws.Cells.Item(33, 9).Value = 13673.667  # I33
$ws.Cells.Item(33, 10).Value = 0  # J33
$ws.Cells.Item(33, 11).Value = 13673.667  # K33
$ws.Cells.Item(33, 12).Value = 0  # L33
$ws.Cells.Item(33, 13).ClearContents()  # M33
$ws.Cells.Item(33, 14).Value = -13337.667  # N33
$ws.Cells.Item(94, 8).Value = 2399.5  # H94
$ws.Cells.Item(94, 9).Value = 1999.1666  # I94
$ws.Cells.Item(94, 11).Value = 1999.1666  # K94
$ws.Cells.Item(94, 13).Value = -1548.1666  # M94

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(47, 8).Value = 27000  # H47
$ws.Cells.Item(47, 9).Value = 17000  # I47
$ws.Cells.Item(47, 11).Value = 17000  # K47
$ws.Cells.Item(47, 13).Value = -16434  # M47

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 299.0625  # H2
$ws.Cells.Item(2, 9).Value = 26  # I2
$ws.Cells.Item(2, 11).Value = 156  # K2
$ws.Cells.Item(2, 13).Value = -43  # M2
$ws.Cells.Item(38, 8).Value = 1115.125  # H38
$ws.Cells.Item(38, 10).Value = 2223.25  # J38
$ws.Cells.Item(38, 12).Value = 6669.75  # L38
$ws.Cells.Item(38, 14).Value = -7363.75  # N38
$ws.Cells.Item(51, 8).Value = 3999.5  # H51
$ws.Cells.Item(51, 9).Value = 3999  # I51
$ws.Cells.Item(51, 11).Value = 11997  # K51
$ws.Cells.Item(51, 13).Value = -11537  # M51
$ws.Cells.Item(88, 8).Value = 0  # H88
$ws.Cells.Item(88, 10).Value = 0  # J88
$ws.Cells.Item(88, 12).ClearContents()  # L88
$ws.Cells.Item(88, 14).Value = 0  # N88
$ws.Cells.Item(91, 8).Value = 0  # H91
$ws.Cells.Item(91, 10).Value = 0  # J91
$ws.Cells.Item(91, 12).ClearContents()  # L91
$ws.Cells.Item(91, 14).Value = 0  # N91
$ws.Cells.Item(132, 8).Value = 901.25  # H132
$ws.Cells.Item(132, 9).Value = 0  # I132
$ws.Cells.Item(132, 10).Value = 901.25  # J132
$ws.Cells.Item(132, 11).Value = 0  # K132
$ws.Cells.Item(132, 12).ClearContents()  # L132
$ws.Cells.Item(132, 13).Value = 8111.25  # M132
$ws.Cells.Item(132, 14).Value = -13171.25  # N132

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(95, 8).Value = 43769.832  # H95
$ws.Cells.Item(95, 10).Value = 43769.832  # J95
$ws.Cells.Item(95, 12).Value = 43769.832  # L95
$ws.Cells.Item(95, 14).Value = -49261.832  # N95
$ws.Cells.Item(122, 8).Value = 7249.75  # H122
$ws.Cells.Item(122, 9).Value = 7666.3335  # I122
$ws.Cells.Item(122, 11).Value = 22999.0005  # K122
$ws.Cells.Item(122, 13).Value = -20549.0005  # M122
$ws.Cells.Item(132, 8).Value = 2025.6666  # H132
$ws.Cells.Item(132, 9).Value = 1638.7142  # I132
$ws.Cells.Item(132, 10).Value = 3380  # J132
$ws.Cells.Item(132, 11).Value = 4916.142599999999  # K132
$ws.Cells.Item(132, 12).Value = 10140  # L132
$ws.Cells.Item(132, 13).Value = -2386.142599999999  # M132
$ws.Cells.Item(132, 14).Value = -15200  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6035.2  # H7
$ws.Cells.Item(7, 9).Value = 4825.3335  # I7
$ws.Cells.Item(7, 10).Value = 7850  # J7
$ws.Cells.Item(7, 11).Value = 4825.3335  # K7
$ws.Cells.Item(7, 12).Value = 7850  # L7
$ws.Cells.Item(7, 13).Value = -4713.3335  # M7
$ws.Cells.Item(7, 14).Value = -8074  # N7
$ws.Cells.Item(16, 8).Value = 60000  # H16
$ws.Cells.Item(16, 9).Value = 0  # I16
$ws.Cells.Item(16, 11).Value = 0  # K16
$ws.Cells.Item(16, 13).ClearContents()  # M16
$ws.Cells.Item(20, 8).Value = 8100  # H20
$ws.Cells.Item(20, 9).Value = 1200  # I20
$ws.Cells.Item(20, 10).Value = 15000  # J20
$ws.Cells.Item(20, 11).Value = 1200  # K20
$ws.Cells.Item(20, 12).Value = 15000  # L20
$ws.Cells.Item(20, 13).Value = -974  # M20
$ws.Cells.Item(20, 14).Value = -15452  # N20
$ws.Cells.Item(46, 8).Value = 3571.4285  # H46
$ws.Cells.Item(46, 9).Value = 1000  # I46
$ws.Cells.Item(46, 10).Value = 10000  # J46
$ws.Cells.Item(46, 11).Value = 1000  # K46
$ws.Cells.Item(46, 12).Value = 10000  # L46
$ws.Cells.Item(46, 13).Value = -812  # M46
$ws.Cells.Item(46, 14).Value = -10376  # N46
$ws.Cells.Item(61, 8).Value = 2000  # H61
$ws.Cells.Item(61, 9).Value = 2000  # I61
$ws.Cells.Item(61, 11).Value = 2000  # K61
$ws.Cells.Item(61, 13).Value = -1798  # M61
$ws.Cells.Item(113, 8).Value = 2000  # H113
$ws.Cells.Item(113, 9).Value = 2000  # I113
$ws.Cells.Item(113, 11).Value = 2000  # K113
$ws.Cells.Item(113, 13).Value = 170  # M113
$ws.Cells.Item(122, 8).Value = 6165.2856  # H122
$ws.Cells.Item(122, 9).Value = 8408.571  # I122
$ws.Cells.Item(122, 11).Value = 25225.713  # K122
$ws.Cells.Item(122, 13).Value = -22775.713  # M122
$ws.Cells.Item(126, 8).Value = 6035.2  # H126
$ws.Cells.Item(126, 9).Value = 4825.3335  # I126
$ws.Cells.Item(126, 10).Value = 7850  # J126
$ws.Cells.Item(126, 11).Value = 14476.0005  # K126
$ws.Cells.Item(126, 12).Value = 23550  # L126
$ws.Cells.Item(126, 13).Value = -12006.0005  # M126
$ws.Cells.Item(126, 14).Value = -28490  # N126
$ws.Cells.Item(132, 8).Value = 3212  # H132
$ws.Cells.Item(132, 9).Value = 3701.6  # I132
$ws.Cells.Item(132, 10).Value = 2600  # J132
$ws.Cells.Item(132, 11).Value = 11104.8  # K132
$ws.Cells.Item(132, 12).Value = 7800  # L132
$ws.Cells.Item(132, 13).Value = -8574.799999999999  # M132
$ws.Cells.Item(132, 14).Value = -12860  # N132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(30, 8).Value = 5499.5  # H30
$ws.Cells.Item(30, 9).Value = 5499.5  # I30
$ws.Cells.Item(30, 10).Value = 0  # J30
$ws.Cells.Item(30, 11).Value = 5499.5  # K30
$ws.Cells.Item(30, 12).Value = 0  # L30
$ws.Cells.Item(30, 13).ClearContents()  # M30
$ws.Cells.Item(30, 14).Value = -5392.5  # N30
$ws.Cells.Item(132, 8).Value = 1165.0358  # H132
$ws.Cells.Item(132, 9).Value = 1165.0358  # I132
$ws.Cells.Item(132, 11).Value = 3495.1074  # K132
$ws.Cells.Item(132, 13).Value = -965.1074000000003  # M132
